# repull data, push all data, mean calculation
# Update the dSF column (F) values for the rows whose underlying data
# was repulled, so that the mean calculation reflects the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = -4
    11 = -8
    19 = -5
    20 = 0
    21 = 0
    37 = -5
    40 = -1
    41 = -2
    42 = -1
    43 = 0
    47 = 0
    48 = -3
    50 = -5
    53 = -3
    58 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
